$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was Cololabis saira / Pacific Saury / Teleost Fish
# becomes Fundulus heteroclitus or majalis / Mummichog or striped killifish / Teleost Fish
$ws.Range("A2").Value = "Fundulus heteroclitus or majalis"
$ws.Range("B2").Value = "Mummichog or striped killifish"
$ws.Range("C2").Value = "Teleost Fish"

# Row 3: was Unassigned / Unassigned / Unassigned
# becomes Cololabis saira / Pacific suary / Teleost Fish
$ws.Range("A3").Value = "Cololabis saira"
$ws.Range("B3").Value = "Pacific suary"
$ws.Range("C3").Value = "Teleost Fish"

# Row 4 (new): Unassigned / Unassigned / Unassigned
$ws.Range("A4").Value = "Unassigned"
$ws.Range("B4").Value = "Unassigned"
$ws.Range("C4").Value = "Unassigned"

# Row 5 (new): Mareca americana / American wigeon / Bird
$ws.Range("A5").Value = "Mareca americana"
$ws.Range("B5").Value = "American wigeon"
$ws.Range("C5").Value = "Bird"

# Row 6 (new): Myrophis vafer / Pacific worm eel / Teleost Fish
$ws.Range("A6").Value = "Myrophis vafer"
$ws.Range("B6").Value = "Pacific worm eel"
$ws.Range("C6").Value = "Teleost Fish"

# Widen columns A:C to fit the longer common names (target stored width 34.5546875;
# closest attainable via the ColumnWidth->stored-width rounding is 34.5)
$ws.Columns.Item(1).ColumnWidth = 33.65
$ws.Columns.Item(2).ColumnWidth = 33.65
$ws.Columns.Item(3).ColumnWidth = 33.65

# Move the active selection to B4, matching the saved workbook state
$ws.Range("B4").Select()
